# Update cryptocurrency price/volume data per the Jan 14 2023 GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "306.03"
Set-TextValue "D3" "32.58"
Set-TextValue "E3" "10.86%"
Set-TextValue "D4" "5.331"
Set-TextValue "E4" "3.72%"
Set-TextValue "D5" "0.07418"
Set-TextValue "E5" "11.31%"
Set-TextValue "D6" "7.768"
Set-TextValue "E6" "5.76%"
Set-TextValue "D7" "3.710"
Set-TextValue "E7" "9.03%"
Set-TextValue "D8" "1.572"
Set-TextValue "E8" "15.75%"
Set-TextValue "D9" "0.9209"
Set-TextValue "E9" "0.38%"
Set-TextValue "D10" "0.01620"
Set-TextValue "E10" "2,413.81%"
Set-TextValue "D11" "0.1671"
Set-TextValue "E11" "6.25%"
Set-TextValue "D12" "0.07395"
Set-TextValue "E12" "12.98%"
Set-TextValue "D13" "0.07980"
Set-TextValue "E13" "4.49%"
Set-TextValue "D14" "0.03118"
Set-TextValue "E14" "6.37%"
Set-TextValue "D15" "0.09819"
Set-TextValue "E15" "9.17%"
Set-TextValue "D16" "0.001521"
Set-TextValue "E16" "-4.01%"
Set-TextValue "D17" "0.04545"
Set-TextValue "E17" "1.53%"
Set-TextValue "D18" "0.006231"
Set-TextValue "E18" "-0.39%"
Set-TextValue "D19" "3.475"
Set-TextValue "E19" "0.53%"
Set-TextValue "D20" "2.241"
Set-TextValue "E20" "0.49%"
Set-TextValue "D21" "0.3273"
Set-TextValue "E21" "1.87%"
Set-TextValue "D22" "0.1321"
Set-TextValue "E22" "0.94%"
Set-TextValue "D23" "4.245"
Set-TextValue "E23" "4.20%"
Set-TextValue "D24" "0.1639"
Set-TextValue "E24" "5.73%"
Set-TextValue "D25" "0.001225"
Set-TextValue "E25" "2.91%"
Set-TextValue "D26" "0.004531"
Set-TextValue "E26" "9.57%"
Set-TextValue "D27" "0.0001168"
Set-TextValue "E27" "-6.45%"
Set-TextValue "D28" "0.0001666"
Set-TextValue "E28" "3.08%"
Set-TextValue "D40" "0.04499"
Set-TextValue "E40" "7.05%"
Set-TextValue "D41" "0.007274"
Set-TextValue "E41" "8.03%"
Set-TextValue "D42" "0.1367"
Set-TextValue "E42" "9.88%"
Set-TextValue "D43" "0.002306"
Set-TextValue "E43" "16.60%"
Set-TextValue "D44" "0.01369"
Set-TextValue "E44" "8.23%"
Set-TextValue "D45" "0.00005971"
Set-TextValue "E45" "5.93%"
Set-TextValue "D47" "0.01300"
Set-TextValue "E47" "-0.40%"
